$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$enDash = [char]0x2013

# Rebuild the whole placeholder body text (existing bullets + the new ones) in one
# shot so every paragraph starts out "clean" (no inherited/explicit indent level).
$fullText = "Lack of Quantum Channel`t" + `
    "`rQuantum Key Distribution use of quantum channel to transport qubits from different users" + `
    "`rAssumption on Quantum Channel" + `
    "`rNo multiple measurements in one circuit " + $enDash + "> every potential key bit is a new experiment" + `
    "`rVery slow" + `
    "`rResource draining" + `
    "`rNoise can appear as "

$tr.Text = $fullText

# Re-apply the outline levels for the paragraphs that need to be indented.
$tr.Paragraphs(2).IndentLevel = 2
$tr.Paragraphs(3).IndentLevel = 2
$tr.Paragraphs(5).IndentLevel = 2
$tr.Paragraphs(6).IndentLevel = 2

# Split the last paragraph into two runs: "Noise can appear as " + "evesdropping".
$tr.InsertAfter("evesdropping") | Out-Null
